$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3743903333333334
$ws.Range("H2").Value = 1.123171
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 23.63579766666667
$ws.Range("N2").Value = 70.907393
$ws.Range("O2").Value = 0.06827844587621175
$ws.Range("P2").Value = 0.06827844587621175
$ws.Range("Q2").Value = 8.849014167022556
$ws.Range("R2").Value = 79.64112750320301
$ws.Range("S2").Value = 0.06827844587621175
$ws.Range("T2").Value = 0.06827844587621175
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.3743903333333334
$ws.Range("H3").Value = 1.123171
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 181.2883913333334
$ws.Range("N3").Value = 543.865174
$ws.Range("O3").Value = 0.5237009467675041
$ws.Range("P3").Value = 0.523700946767504
$ws.Range("Q3").Value = 67.87262126075046
$ws.Range("R3").Value = 610.8535913467541
$ws.Range("S3").Value = 0.5237009467675041
$ws.Range("T3").Value = 0.523700946767504
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.3743903333333334
$ws.Range("H4").Value = 1.123171
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 111.1005463333333
$ws.Range("N4").Value = 333.301639
$ws.Range("O4").Value = 0.3209442197221123
$ws.Range("P4").Value = 0.3209442197221123
$ws.Range("Q4").Value = 41.59497057525211
$ws.Range("R4").Value = 374.354735177269
$ws.Range("S4").Value = 0.3209442197221123
$ws.Range("T4").Value = 0.3209442197221123
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.3743903333333334
$ws.Range("H5").Value = 1.123171
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 30.14303933333333
$ws.Range("N5").Value = 90.42911799999999
$ws.Range("O5").Value = 0.08707638763417187
$ws.Range("P5").Value = 0.08707638763417187
$ws.Range("Q5").Value = 11.28526254368644
$ws.Range("R5").Value = 101.567362893178
$ws.Range("S5").Value = 0.08707638763417187
$ws.Range("T5").Value = 0.08707638763417187
